$d = $word.ActiveDocument

$replacements = @(
    @("2024-05-05 Sunday", "2024-05-06 Monday"),
    @("45×55=", "49×19="),
    @("41×14=", "54×91="),
    @("40×81=", "83×90="),
    @("89×29=", "72×31="),
    @("50×11=", "11×60="),
    @("63×46=", "88×62="),
    @("37×52=", "71×75="),
    @("99×29=", "85×84="),
    @("99×75=", "65×36="),
    @("16×65=", "97×78="),
    @("52×99=", "51×13="),
    @("11×97=", "64×87="),
    @("97×27=", "55×71="),
    @("37×93=", "70×32="),
    @("64×50=", "52×70="),
    @("28×59=", "53×92="),
    @("45×46=", "50×99="),
    @("35×74=", "43×77="),
    @("21×52=", "20×64="),
    @("35×96=", "55×19="),
    @("57×38=", "50×46="),
    @("52×36=", "87×60="),
    @("72×60=", "39×36="),
    @("65×18=", "15×87="),
    @("76×55=", "72×48=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
